$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels for the new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style from an existing header cell (e.g. AC1) onto the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Fill the team record values for each data row (rows 2-50)
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 73   # AD = column 30
    $ws.Cells.Item($r, 31).Value = 89   # AE = column 31
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32
}
